$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row at row 12 (push old rows 12-15 down to 13-16), formatted like row 11 ---
$ws.Range("A11:E11").Copy()
$ws.Range("A12:E12").Insert(-4121)
$excel.CutCopyMode = $false

# The simulated insert can drop the right-hand border on C12; re-copy the exact
# per-cell format from row 11 so every cell in the new row matches row 11 precisely.
$ws.Range("A11").Copy()
$ws.Range("A12").PasteSpecial(-4122)
$ws.Range("B11").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("C11").Copy()
$ws.Range("C12").PasteSpecial(-4122)
$ws.Range("D11").Copy()
$ws.Range("D12").PasteSpecial(-4122)
$ws.Range("E11").Copy()
$ws.Range("E12").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height for the new row should match row 11's (18pt, custom height)
$ws.Rows.Item(12).RowHeight = $ws.Rows.Item(11).RowHeight

# --- Fill in the new "Half Year to 31 Dec 2016" data row ---
$ws.Range("A12").Value = "Half Year to"
$ws.Range("B12").Value = 42735
$ws.Range("C12").Value = 13.6
$ws.Range("D12").Value = 10.4
$ws.Range("E12").Value = 3.1

# --- Update the recalculated summary rows (now pushed down to 13 & 14) ---
# Row 13: "Annualised Performance"
$ws.Range("C13").Value = 16.2
$ws.Range("D13").Value = 5.5
$ws.Range("E13").Value = 10.7

# Row 14: "Cumulative Performance"
$ws.Range("C14").Value = 385
$ws.Range("D14").Value = 76
$ws.Range("E14").Value = 309

# --- Reset the sheet view: scroll back to the top and select A1 (was topLeftCell=A2 / F11) ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("A1").Select()
